# Update cryptocurrency price/volume data per the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.116.82"
$ws.Range("E2").Value = "  -0.66%  "
$ws.Range("D3").Value = "2.445.53"
$ws.Range("E3").Value = "  +0.37%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.67"
$ws.Range("E5").Value = "  +1.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.58"
$ws.Range("E6").Value = "  -1.07%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.533"
$ws.Range("E8").Value = "  +0.29%  "
$ws.Range("D9").Value = "2.440.32"
$ws.Range("E9").Value = "  +0.34%  "
$ws.Range("E10").Value = "  +1.45%  "
$ws.Range("E12").Value = "  -0.20%  "
$ws.Range("E13").Value = "  -2.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.37"
$ws.Range("E14").Value = "  -0.83%  "
$ws.Range("E15").Value = "  +0.30%  "
$ws.Range("D16").Value = "2.876.44"
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("D17").Value = "62.039.77"
$ws.Range("D18").Value = "2.436.39"
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.73"
$ws.Range("E19").Value = "  -3.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.24"
$ws.Range("E20").Value = "  +1.39%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "326.19"
$ws.Range("E21").Value = "  -0.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.10"
$ws.Range("E22").Value = "  -1.19%  "
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("E24").Value = "  -6.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.57"
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.11"
$ws.Range("E26").Value = "  +0.44%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "599.68"
$ws.Range("E27").Value = "  -5.20%  "
$ws.Range("D28").Value = "0.0₃0964"
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("D29").Value = "2.564.50"
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("E30").Value = "  -0.45%  "
$ws.Range("E31").Value = "  -1.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.97"
$ws.Range("E32").Value = "  -1.09%  "
$ws.Range("E33").Value = "  +0.49%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.88"
$ws.Range("E35").Value = "  -2.31%  "
$ws.Range("E36").Value = "  +0.25%  "
$ws.Range("E37").Value = "  -1.81%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.376"
$ws.Range("E38").Value = "  +0.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "152.77"
$ws.Range("E39").Value = "  +3.96%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.39"
$ws.Range("E40").Value = "  -0.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.28"
$ws.Range("E41").Value = "  +0.65%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "43.14"
$ws.Range("E42").Value = "  +2.02%  "
$ws.Range("E43").Value = "  -1.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.50"
$ws.Range("E45").Value = "  +0.70%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "142.01"
$ws.Range("E46").Value = "  -2.41%  "
$ws.Range("E47").Value = "  -2.62%  "
$ws.Range("E48").Value = "  +17.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.601"
$ws.Range("E49").Value = "  +0.69%  "
$ws.Range("E50").Value = "  -0.94%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.80"
$ws.Range("E51").Value = "  +0.13%  "
